# bd/ruta.xlsx — "Agregando app con modificaciones"
#
# The ACHS Osorno stop (row 3) was removed from the route, which shifts
# every following stop up by one row. Two brand-new stops (Maria José
# Rodriguez / Isaias Beroiza Mora, both in Calbuco) were appended at the
# end of the now-shorter list, landing on rows 18-19. The trailing blank
# row 138 disappears as a consequence (dimension becomes L137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the ACHS Osorno row entirely; Excel shifts rows 4:138 up to 3:137
# (this alone also fixes row 64's stray formatting and drops the blank
# row that used to be 138).
$ws.Rows.Item(3).Delete()

# --- New stop (row 18): Maria José Rodriguez ---------------------------
# A18/B18 are numeric columns even though the column's number format is
# Text ("@"); flip to General while assigning so Excel stores a real
# number instead of coercing it to a text string, then restore the
# column's Text format.
$ws.Range("A18").NumberFormat = "General"
$ws.Range("A18").Value = 20250318
$ws.Range("A18").NumberFormat = "@"

$ws.Range("B18").NumberFormat = "General"
$ws.Range("B18").Value = 16
$ws.Range("B18").NumberFormat = "@"

$ws.Range("C18").Value = "30"
$ws.Range("D18").Value = "17.673.326-8"
$ws.Range("E18").Value = "Maria José Rodriguez"
$ws.Range("F18").Value = "Colaco s/n km 3, parcela 9"
$ws.Range("G18").Value = "Calbuco"
$ws.Range("H18").Value = "972861950"
$ws.Range("I18").Value = "Cliente test"
$ws.Range("J18").Value = "1002"

# --- New stop (row 19): Isaias Beroiza Mora -----------------------------
$ws.Range("A19").NumberFormat = "General"
$ws.Range("A19").Value = 20250318
$ws.Range("A19").NumberFormat = "@"

$ws.Range("B19").NumberFormat = "General"
$ws.Range("B19").Value = 17
$ws.Range("B19").NumberFormat = "@"

$ws.Range("C19").Value = "30"
$ws.Range("D19").Value = "16.742.249-7"
$ws.Range("E19").Value = "Isaias Beroiza Mora"
$ws.Range("F19").Value = "Colaco s/n km 3, parcela 9"
$ws.Range("G19").Value = "Calbuco"
$ws.Range("H19").Value = "988809704"
$ws.Range("I19").Value = "Cliente test"
$ws.Range("J19").Value = "1001"
